$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "NSE:BAJAJ-AUTO"
$ws.Range("C2").Value = "NSE:AAATECH"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "NSE:BPCL"
$ws.Range("F2").Value = "NSE:BAJAJ-AUTO"

# Row 3
$ws.Range("B3").Value = "NSE:BINANIIND"
$ws.Range("C3").Value = "NSE:ADVANIHOTR"
$ws.Range("E3").Value = "NSE:LICHSGFIN"

# Row 4
$ws.Range("B4").Value = "NSE:BRITANNIA"
$ws.Range("C4").Value = "NSE:AGSTRA"

# Row 5
$ws.Range("B5").Value = "NSE:CAMPUS"
$ws.Range("C5").Value = "NSE:CERA"

# Row 6
$ws.Range("B6").Value = "NSE:DELTAMAGNT"
$ws.Range("C6").Value = "NSE:GINNIFILA"

# Row 7
$ws.Range("C7").Value = "NSE:KPRMILL"

# Row 8
$ws.Range("B8").Value = "NSE:DIXON"
$ws.Range("C8").Value = "NSE:LAL"

# Row 9
$ws.Range("B9").Value = "NSE:DTIL"
$ws.Range("C9").Value = "NSE:MANUGRAPH"

# Row 10
$ws.Range("B10").Value = "NSE:GULPOLY"
$ws.Range("C10").Value = "NSE:MEGASOFT"

# Row 11
$ws.Range("B11").Value = "NSE:MEDICAMEQ"
$ws.Range("C11").Value = "NSE:NITIRAJ"

# Row 12
$ws.Range("C12").Value = "NSE:PFC"

# Row 13
$ws.Range("C13").Value = "NSE:PSB"

# Row 14
$ws.Range("C14").Value = "NSE:QUICKHEAL"

# Row 15
$ws.Range("C15").Value = "NSE:RECLTD"

# Remove rows 16-23 (no longer present in the updated data set)
$ws.Range("A16:F23").EntireRow.Delete()
